$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149, shifting rows 149:242 down to 150:243
$ws.Rows.Item(149).Insert()

# Populate the new row 149 with data
$ws.Cells.Item(149, 1).Value = 9
$ws.Cells.Item(149, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(149, 3).Value = "Metropolitana"
$ws.Cells.Item(149, 4).Value = 44609
$ws.Cells.Item(149, 5).Value = 13
$ws.Cells.Item(149, 6).Value = 100112001
$ws.Cells.Item(149, 7).Value = "Berenjena"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 97
$ws.Cells.Item(149, 11).Value = 10000
$ws.Cells.Item(149, 12).Value = 12000
$ws.Cells.Item(149, 13).Value = 10990
$ws.Cells.Item(149, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(149, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(149, 16).Value = 183
$ws.Cells.Item(149, 17).Value = 60
$ws.Cells.Item(149, 18).Value = "Hortaliza"
